$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 15.81542716379291
$ws.Cells.Item(2, 4).Value = 4.608468162201127
$ws.Cells.Item(2, 5).Value = 20.52734289458135
$ws.Cells.Item(2, 6).Value = 22.23093670878155
$ws.Cells.Item(2, 7).Value = 25.22786576219513
$ws.Cells.Item(2, 8).Value = 13.41310485122486
$ws.Cells.Item(2, 11).Value = 10.94344334586449
$ws.Cells.Item(2, 12).Value = 8.047874645982276
$ws.Cells.Item(2, 13).Value = 14.68701805929012
$ws.Cells.Item(2, 14).Value = 20.14463931274576
$ws.Cells.Item(2, 15).Value = 19.94254533248236

$ws.Cells.Item(3, 2).Value = 15.73861400343886
$ws.Cells.Item(3, 4).Value = 4.530428872322533
$ws.Cells.Item(3, 5).Value = 20.59599441593353
$ws.Cells.Item(3, 6).Value = 22.23415789061444
$ws.Cells.Item(3, 7).Value = 25.22703743072563
$ws.Cells.Item(3, 8).Value = 13.44830304674604
$ws.Cells.Item(3, 11).Value = 10.71212393008241
$ws.Cells.Item(3, 12).Value = 8.02941573471257
$ws.Cells.Item(3, 13).Value = 14.66991148258432
$ws.Cells.Item(3, 14).Value = 20.2023779973051
$ws.Cells.Item(3, 15).Value = 19.99014779165025

$ws.Cells.Item(4, 2).Value = 15.69422573477628
$ws.Cells.Item(4, 4).Value = 4.481106919419823
$ws.Cells.Item(4, 5).Value = 20.64057310338942
$ws.Cells.Item(4, 6).Value = 22.24170720810843
$ws.Cells.Item(4, 7).Value = 25.23426941089593
$ws.Cells.Item(4, 8).Value = 13.47185678207637
$ws.Cells.Item(4, 11).Value = 10.56599773511888
$ws.Cells.Item(4, 12).Value = 8.019123767476959
$ws.Cells.Item(4, 13).Value = 14.66142799142747
$ws.Cells.Item(4, 14).Value = 20.23957544135426
$ws.Cells.Item(4, 15).Value = 20.02334745643375

$ws.Cells.Item(5, 2).Value = 15.67685017611077
$ws.Cells.Item(5, 4).Value = 4.460666973944446
$ws.Cells.Item(5, 5).Value = 20.65935080185002
$ws.Cells.Item(5, 6).Value = 22.24618547218956
$ws.Cells.Item(5, 7).Value = 25.23916239248683
$ws.Cells.Item(5, 8).Value = 13.48194354493926
$ws.Cells.Item(5, 11).Value = 10.50546880680027
$ws.Cells.Item(5, 12).Value = 8.015194624801723
$ws.Cells.Item(5, 13).Value = 14.65848219439608
$ws.Cells.Item(5, 14).Value = 20.25517386533287
$ws.Cells.Item(5, 15).Value = 20.03787406376721

$ws.Cells.Item(6, 2).Value = 15.67400846693352
$ws.Cells.Item(6, 4).Value = 4.457252728146177
$ws.Cells.Item(6, 5).Value = 20.66250579975698
$ws.Cells.Item(6, 6).Value = 22.24701376181039
$ws.Cells.Item(6, 7).Value = 25.24009234995719
$ws.Cells.Item(6, 8).Value = 13.48364794712746
$ws.Cells.Item(6, 11).Value = 10.49536022297135
$ws.Cells.Item(6, 12).Value = 8.014558275911916
$ws.Cells.Item(6, 13).Value = 14.65802402057967
$ws.Cells.Item(6, 14).Value = 20.25779059200924
$ws.Cells.Item(6, 15).Value = 20.04034640278835

$ws.Cells.Item(7, 2).Value = 15.69398849559682
$ws.Cells.Item(7, 4).Value = 4.480832622251348
$ws.Cells.Item(7, 5).Value = 20.64082386810882
$ws.Cells.Item(7, 6).Value = 22.24176192697159
$ws.Cells.Item(7, 7).Value = 25.23432752235889
$ws.Cells.Item(7, 8).Value = 13.47199083792766
$ws.Cells.Item(7, 11).Value = 10.565185328194
$ws.Cells.Item(7, 12).Value = 8.019069701312953
$ws.Cells.Item(7, 13).Value = 14.66138618893934
$ws.Cells.Item(7, 14).Value = 20.23978402334837
$ws.Cells.Item(7, 15).Value = 20.02353933030648

$ws.Cells.Item(8, 2).Value = 15.78837532770557
$ws.Cells.Item(8, 4).Value = 4.581859304261196
$ws.Cells.Item(8, 5).Value = 20.55051118193479
$ws.Cells.Item(8, 6).Value = 22.23089176738666
$ws.Cells.Item(8, 7).Value = 25.22597372107731
$ws.Cells.Item(8, 8).Value = 13.42483822203654
$ws.Cells.Item(8, 11).Value = 10.86456257716164
$ws.Cells.Item(8, 12).Value = 8.041295582017698
$ws.Cells.Item(8, 13).Value = 14.68070253820644
$ws.Cells.Item(8, 14).Value = 20.16418602668024
$ws.Cells.Item(8, 15).Value = 19.95813367756529

$ws.Cells.Item(9, 2).Value = 15.99478238763502
$ws.Cells.Item(9, 4).Value = 4.768297147498711
$ws.Cells.Item(9, 5).Value = 20.39259714875124
$ws.Cells.Item(9, 6).Value = 22.25371040317478
$ws.Cells.Item(9, 7).Value = 25.27096753071489
$ws.Cells.Item(9, 8).Value = 13.34777476731676
$ws.Cells.Item(9, 11).Value = 11.41716057749484
$ws.Cells.Item(9, 12).Value = 8.093011474396834
$ws.Cells.Item(9, 13).Value = 14.73445329292732
$ws.Cells.Item(9, 14).Value = 20.02973399332963
$ws.Cells.Item(9, 15).Value = 19.8614352015176

$ws.Cells.Item(10, 2).Value = 16.15842291479923
$ws.Cells.Item(10, 4).Value = 4.897497997597998
$ws.Cells.Item(10, 5).Value = 20.28818722383339
$ws.Cells.Item(10, 6).Value = 22.2972325827756
$ws.Cells.Item(10, 7).Value = 25.34128563556391
$ws.Cells.Item(10, 8).Value = 13.30053839218741
$ws.Cells.Item(10, 11).Value = 11.79967568671665
$ws.Cells.Item(10, 12).Value = 8.135771765637042
$ws.Cells.Item(10, 13).Value = 14.78339601447744
$ws.Cells.Item(10, 14).Value = 19.93928332312156
$ws.Cells.Item(10, 15).Value = 19.80969186497671

$ws.Cells.Item(11, 2).Value = 16.23524306830631
$ws.Cells.Item(11, 4).Value = 4.954459333498538
$ws.Cells.Item(11, 5).Value = 20.24319064936543
$ws.Cells.Item(11, 6).Value = 22.32279197865148
$ws.Cells.Item(11, 7).Value = 25.38129947817476
$ws.Cells.Item(11, 8).Value = 13.28108474154479
$ws.Cells.Item(11, 11).Value = 11.96813008221014
$ws.Cells.Item(11, 12).Value = 8.15621456805774
$ws.Cells.Item(11, 13).Value = 14.80766080235702
$ws.Cells.Item(11, 14).Value = 19.8999272662915
$ws.Cells.Item(11, 15).Value = 19.79035315476687

$ws.Cells.Item(12, 2).Value = 16.26465480452833
$ws.Cells.Item(12, 4).Value = 4.975758553729015
$ws.Cells.Item(12, 5).Value = 20.22650973133192
$ws.Cells.Item(12, 6).Value = 22.33329331650925
$ws.Cells.Item(12, 7).Value = 25.39759748996058
$ws.Cells.Item(12, 8).Value = 13.27401054355173
$ws.Cells.Item(12, 11).Value = 12.03108571509962
$ws.Cells.Item(12, 12).Value = 8.164094163776937
$ws.Cells.Item(12, 13).Value = 14.81713197935756
$ws.Cells.Item(12, 14).Value = 19.88528037647432
$ws.Cells.Item(12, 15).Value = 19.78363446450011

$ws.Cells.Item(13, 2).Value = 16.258306489225
$ws.Cells.Item(13, 4).Value = 4.971183584135576
$ws.Cells.Item(13, 5).Value = 20.2300863484498
$ws.Cells.Item(13, 6).Value = 22.33099518551005
$ws.Cells.Item(13, 7).Value = 25.39403662098952
$ws.Cells.Item(13, 8).Value = 13.2755210919329
$ws.Cells.Item(13, 11).Value = 12.01756474315606
$ws.Cells.Item(13, 12).Value = 8.162391066354846
$ws.Cells.Item(13, 13).Value = 14.81507970488838
$ws.Cells.Item(13, 14).Value = 19.88842345983922
$ws.Cells.Item(13, 15).Value = 19.7850545642271

$ws.Cells.Item(14, 2).Value = 16.23765646331318
$ws.Cells.Item(14, 4).Value = 4.956217116217954
$ws.Cells.Item(14, 5).Value = 20.24181112593924
$ws.Cells.Item(14, 6).Value = 22.32363948176456
$ws.Cells.Item(14, 7).Value = 25.38261741937048
$ws.Cells.Item(14, 8).Value = 13.28049688179175
$ws.Cells.Item(14, 11).Value = 11.97332639754741
$ws.Cells.Item(14, 12).Value = 8.156860078694308
$ws.Cells.Item(14, 13).Value = 14.80843436725993
$ws.Cells.Item(14, 14).Value = 19.89871712536524
$ws.Cells.Item(14, 15).Value = 19.78978828725752

$ws.Cells.Item(15, 2).Value = 16.22504899174537
$ws.Cells.Item(15, 4).Value = 4.947014166230991
$ws.Cells.Item(15, 5).Value = 20.24903951790016
$ws.Cells.Item(15, 6).Value = 22.3192408313631
$ws.Cells.Item(15, 7).Value = 25.37577174563057
$ws.Cells.Item(15, 8).Value = 13.28358278360174
$ws.Cells.Item(15, 11).Value = 11.94611945662725
$ws.Cells.Item(15, 12).Value = 8.153490080627272
$ws.Cells.Item(15, 13).Value = 14.80440055874969
$ws.Cells.Item(15, 14).Value = 19.90505565154118
$ws.Cells.Item(15, 15).Value = 19.79276656123776

$ws.Cells.Item(16, 2).Value = 16.15344867132837
$ws.Cells.Item(16, 4).Value = 4.893738082762678
$ws.Cells.Item(16, 5).Value = 20.29117806037413
$ws.Cells.Item(16, 6).Value = 22.29567766031976
$ws.Cells.Item(16, 7).Value = 25.33883145027952
$ws.Cells.Item(16, 8).Value = 13.3018506660999
$ws.Cells.Item(16, 11).Value = 11.78855191977279
$ws.Cells.Item(16, 12).Value = 8.134455369827776
$ws.Cells.Item(16, 13).Value = 14.78185008354636
$ws.Cells.Item(16, 14).Value = 19.94189127035013
$ws.Cells.Item(16, 15).Value = 19.81104025420163

$ws.Cells.Item(17, 2).Value = 16.11011914473951
$ws.Cells.Item(17, 4).Value = 4.860583232538133
$ws.Cells.Item(17, 5).Value = 20.31766815382219
$ws.Cells.Item(17, 6).Value = 22.28269404887292
$ws.Cells.Item(17, 7).Value = 25.31821987191208
$ws.Cells.Item(17, 8).Value = 13.31357843864106
$ws.Cells.Item(17, 11).Value = 11.69044057095779
$ws.Cells.Item(17, 12).Value = 8.123029094493086
$ws.Cells.Item(17, 13).Value = 14.76852496722384
$ws.Cells.Item(17, 14).Value = 19.96494653476234
$ws.Cells.Item(17, 15).Value = 19.82332665030938

$ws.Cells.Item(18, 2).Value = 16.08542235799209
$ws.Cells.Item(18, 4).Value = 4.841343378127913
$ws.Cells.Item(18, 5).Value = 20.33313993928507
$ws.Cells.Item(18, 6).Value = 22.27576887505314
$ws.Cells.Item(18, 7).Value = 25.30712067277008
$ws.Cells.Item(18, 8).Value = 13.32051543696397
$ws.Cells.Item(18, 11).Value = 11.6334887427173
$ws.Cells.Item(18, 12).Value = 8.116550532670919
$ws.Cells.Item(18, 13).Value = 14.76104932704552
$ws.Cells.Item(18, 14).Value = 19.9783758973934
$ws.Cells.Item(18, 15).Value = 19.8307887241352

$ws.Cells.Item(19, 2).Value = 16.07709974727264
$ws.Cells.Item(19, 4).Value = 4.834800204540621
$ws.Cells.Item(19, 5).Value = 20.33841888408605
$ws.Cells.Item(19, 6).Value = 22.27351749567746
$ws.Cells.Item(19, 7).Value = 25.30349275547445
$ws.Cells.Item(19, 8).Value = 13.32289707646883
$ws.Cells.Item(19, 11).Value = 11.61411752459931
$ws.Cells.Item(19, 12).Value = 8.114373190415341
$ws.Cells.Item(19, 13).Value = 14.75855074785119
$ws.Cells.Item(19, 14).Value = 19.98295183597097
$ws.Cells.Item(19, 15).Value = 19.83338312078961

$ws.Cells.Item(20, 2).Value = 16.11470848535659
$ws.Cells.Item(20, 4).Value = 4.864130298987187
$ws.Cells.Item(20, 5).Value = 20.31482388545285
$ws.Cells.Item(20, 6).Value = 22.28402005320432
$ws.Cells.Item(20, 7).Value = 25.32033581886202
$ws.Cells.Item(20, 8).Value = 13.31231017879234
$ws.Cells.Item(20, 11).Value = 11.70093885784482
$ws.Cells.Item(20, 12).Value = 8.124235790574161
$ws.Cells.Item(20, 13).Value = 14.76992396344655
$ws.Cells.Item(20, 14).Value = 19.96247482240301
$ws.Cells.Item(20, 15).Value = 19.82197782946815

$ws.Cells.Item(21, 2).Value = 16.24371331409371
$ws.Cells.Item(21, 4).Value = 4.960620557721969
$ws.Cells.Item(21, 5).Value = 20.23835755961666
$ws.Cells.Item(21, 6).Value = 22.32577775777756
$ws.Cells.Item(21, 7).Value = 25.38594049890498
$ws.Cells.Item(21, 8).Value = 13.27902743457735
$ws.Cells.Item(21, 11).Value = 11.98634318854244
$ws.Cells.Item(21, 12).Value = 8.158480942868763
$ws.Cells.Item(21, 13).Value = 14.81037863497042
$ws.Cells.Item(21, 14).Value = 19.89568667621123
$ws.Cells.Item(21, 15).Value = 19.78838146965519

$ws.Cells.Item(22, 2).Value = 16.32989016902734
$ws.Cells.Item(22, 4).Value = 5.022099822845825
$ws.Cells.Item(22, 5).Value = 20.19047034812551
$ws.Cells.Item(22, 6).Value = 22.35785991652175
$ws.Cells.Item(22, 7).Value = 25.43548953600788
$ws.Cells.Item(22, 8).Value = 13.2589799747754
$ws.Cells.Item(22, 11).Value = 12.16799329278732
$ws.Cells.Item(22, 12).Value = 8.181666508785773
$ws.Cells.Item(22, 13).Value = 14.83846296245853
$ws.Cells.Item(22, 14).Value = 19.85353074210791
$ws.Cells.Item(22, 15).Value = 19.76994768470672

$ws.Cells.Item(23, 2).Value = 16.28373224114937
$ws.Cells.Item(23, 4).Value = 4.989435183516601
$ws.Cells.Item(23, 5).Value = 20.21583799588141
$ws.Cells.Item(23, 6).Value = 22.34030085304559
$ws.Cells.Item(23, 7).Value = 25.40843695491108
$ws.Cells.Item(23, 8).Value = 13.26952372490585
$ws.Cells.Item(23, 11).Value = 12.07150067030535
$ws.Cells.Item(23, 12).Value = 8.16921974440282
$ws.Cells.Item(23, 13).Value = 14.82332505324619
$ws.Cells.Item(23, 14).Value = 19.87589380634421
$ws.Cells.Item(23, 15).Value = 19.77946361684669

$ws.Cells.Item(24, 2).Value = 16.11263297652759
$ws.Cells.Item(24, 4).Value = 4.862527226051174
$ws.Cells.Item(24, 5).Value = 20.31610902443149
$ws.Cells.Item(24, 6).Value = 22.28341888641606
$ws.Cells.Item(24, 7).Value = 25.31937686035016
$ws.Cells.Item(24, 8).Value = 13.31288295306785
$ws.Cells.Item(24, 11).Value = 11.69619428124687
$ws.Cells.Item(24, 12).Value = 8.12368996089382
$ws.Cells.Item(24, 13).Value = 14.76929090016612
$ws.Cells.Item(24, 14).Value = 19.9635917396606
$ws.Cells.Item(24, 15).Value = 19.82258639023991

$ws.Cells.Item(25, 2).Value = 15.93676216760382
$ws.Cells.Item(25, 4).Value = 4.719180880022554
$ws.Cells.Item(25, 5).Value = 20.43327199128807
$ws.Cells.Item(25, 6).Value = 22.2428241911074
$ws.Cells.Item(25, 7).Value = 25.25223445313591
$ws.Cells.Item(25, 8).Value = 13.36697425931284
$ws.Cells.Item(25, 11).Value = 11.27162438794591
$ws.Cells.Item(25, 12).Value = 8.078169599990423
$ws.Cells.Item(25, 13).Value = 14.71823516481356
$ws.Cells.Item(25, 14).Value = 20.06463813069656
$ws.Cells.Item(25, 15).Value = 19.88420893229395
